$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the "Trimestre" column as plain text (not auto-converted to a date value)
$ws.Range("C2:C9").NumberFormat = "@"

# Row 2: Amapá -> Bahia, date -> 01/01/2024, value -> 14
$ws.Range("A2").Value = "Bahia"
$ws.Range("C2").Value = "01/01/2024"
$ws.Range("D2").Value = 14

# Row 3: Bahia -> Pernambuco, date -> 01/01/2024, value -> 12.4
$ws.Range("A3").Value = "Pernambuco"
$ws.Range("C3").Value = "01/01/2024"
$ws.Range("D3").Value = 12.4

# Row 4: Pernambuco -> Amapá, date -> 01/01/2024, value -> 10.9
$ws.Range("A4").Value = "Amapá"
$ws.Range("C4").Value = "01/01/2024"
$ws.Range("D4").Value = 10.9

# Row 5: Sergipe -> Rio de Janeiro, date -> 01/01/2024, value -> 10.3
$ws.Range("A5").Value = "Rio de Janeiro"
$ws.Range("C5").Value = "01/01/2024"
$ws.Range("D5").Value = 10.3

# Row 6: Piauí stays, date -> 01/01/2024, value -> 10
$ws.Range("C6").Value = "01/01/2024"
$ws.Range("D6").Value = 10

# Row 7: Rio de Janeiro -> Sergipe, date -> 01/01/2024, value unchanged (10)
$ws.Range("A7").Value = "Sergipe"
$ws.Range("C7").Value = "01/01/2024"

# Row 8: Nordeste stays, date -> 01/01/2024, value -> 11.1
$ws.Range("C8").Value = "01/01/2024"
$ws.Range("D8").Value = 11.1

# Row 9: Brasil stays, date -> 01/01/2024, value -> 7.9
$ws.Range("C9").Value = "01/01/2024"
$ws.Range("D9").Value = 7.9
